# Insert a new column before column A to make room for the "ID" column,
# then populate it with the header "ID" and per-row identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns (A:E) one position to the right (-> B:F).
$ws.Columns("A:A").Insert()

# Header for the new ID column, matching the style of the other headers
# (bold font, centered/top alignment, thin border) by copying the format
# from the neighboring header cell.
$ws.Range("A1").Value2 = "ID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row identifiers for the new ID column (rows 2-25).
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $ids[$i]
}
